$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix missing TIME OUT value on 04-26-2015 (row 6) - same as TIME IN
$ws.Range("D6").Value = "09:39:30"

# Fill in the Official Business (OB Meeting) start/end times that were left
# blank for the four OB Meeting days (rows 7-10):
#   K/L = OFFICIAL BUSINESS DEPARTURE / TIME START -> 08:30:00
#   M/N = OFFICIAL BUSINESS TIME END / ARRIVAL     -> 18:30:00
$obRows = 7, 8, 9, 10
foreach ($r in $obRows) {
    $ws.Range("K$r").Value = "08:30:00"
    $ws.Range("L$r").Value = "08:30:00"
    $ws.Range("M$r").Value = "18:30:00"
    $ws.Range("N$r").Value = "18:30:00"
}
